# ----------------------------------------------------------------------
# Checkpoint 10 de Maio.pptx - applies the edits described by the diff:
#   1. Slide 1 ("TextBox 21"): fill the empty 2nd paragraph with
#      "Unidade Curricular: Inteligência Artificial" (4 runs), which
#      also makes spAutoFit shrink the shape's height.
#   2-6. Slides 2-6 (title "TextBox 7" shapes): bump every run's font
#      size from 24pt to 28pt, which also makes spAutoFit grow the
#      shape's height.
# ----------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Slide 1, shape 5 ("TextBox 21") - title textbox on the cover slide
# ------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$titleBox = $slide1.Shapes.Item(5)
$tr1 = $titleBox.TextFrame.TextRange

# Paragraph 2 is currently empty (just a paragraph mark) - insert the
# new "Unidade Curricular: Inteligência Artificial" text into it.
$para2 = $tr1.Paragraphs(2, 1)
$para2.InsertAfter("Unidade Curricular: Inteligência Artificial")

# Re-fetch the whole range text so the new characters can be addressed,
# then apply the Optima latin typeface to each of the four runs.
$tr1b = $titleBox.TextFrame.TextRange
$run1 = $tr1b.Characters(37, 7)   # "Unidade"
$run1.Font.Name = "Optima"
$run2 = $tr1b.Characters(44, 13)  # " Curricular: "
$run2.Font.Name = "Optima"
$run3 = $tr1b.Characters(57, 12)  # "Inteligência"
$run3.Font.Name = "Optima"
$run4 = $tr1b.Characters(69, 11)  # " Artificial"
$run4.Font.Name = "Optima"

# ------------------------------------------------------------------
# 2-6) Slides 2-6, the title shape ("TextBox 7") - bump 24pt -> 28pt
# ------------------------------------------------------------------
$titleSlideIdx = @(2, 3, 4, 5, 6)

foreach ($idx in $titleSlideIdx) {
    $slide = $p.Slides.Item($idx)
    $shape = $slide.Shapes.Item(6)
    $tr = $shape.TextFrame.TextRange

    # Bump every existing run in the title to 28pt.
    $tr.Font.Size = 28

    # The shape uses <a:spAutoFit/>, so its cached height only gets
    # recomputed once the text actually changes. Force a reflow by
    # appending then deleting a dummy character.
    $tr.Text = $tr.Text + "X"
    $tail = $tr.Characters($tr.Text.Length, 1)
    $tail.Text = ""
}
